$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "CasesTab" row's query (B2) is being updated to drop the trailing
# `Cohort` column from its RETURN clause (the Cypher query text literally
# loses the `coalesce(co.cohort_description, '') AS `Cohort`` line).
$caseQueryNoCohort = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
MATCH (samp:sample)-->(c) 
WHERE samp.summarized_sample_type IN ["Primary Malignant Tumor Tissue"] 
OPTIONAL MATCH (co:cohort)<-[*]-(c)
  WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $caseQueryNoCohort

# The text in B2 now wraps over one fewer line (12 instead of 13), and all
# three data rows settle on the same wrapped height once re-measured.
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 244.8
$ws.Rows.Item(4).RowHeight = 244.8

# Selection moves from the old FilesTab row (B4) up to the CasesTab row
# (B2) that was just edited, so that's what is left active/selected and
# in view (no more scrolled-down topLeftCell).
$ws.Range("B2").Select()
